# Fix terminology mappings on the "Properties" sheet of the
# HeartRateVariabilityCS CodeSystem workbook:
#   - loinc-equivalent property: Uri now points at the
#     concept-property-definitions-cs CodeSystem fragment, and its Type
#     is corrected from "code" to "string".
#   - status property is renamed to "assignment-status", its Uri is
#     updated to the concept-property-definitions-cs CodeSystem fragment,
#     and its Description is reworded from "Status of LOINC code
#     assignment" to "Status of terminology code assignment".

$wb = $excel.ActiveWorkbook
$props = $wb.Worksheets.Item("Properties")

# Row 2: loinc-equivalent property (Code | Uri | Description | Type)
$props.Range("B2").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/CodeSystem/concept-property-definitions-cs#loinc-equivalent"
$props.Range("D2").Value = "string"

# Row 3: status -> assignment-status property (Code | Uri | Description | Type)
$props.Range("A3").Value = "assignment-status"
$props.Range("B3").Value = "https://2rdoc.pt/ig/ios-lifestyle-medicine/CodeSystem/concept-property-definitions-cs#assignment-status"
$props.Range("C3").Value = "Status of terminology code assignment"
